# Rename the single worksheet from "مثال على الاخراج" to "البيانات".
# Excel automatically repoints unqualified/qualified references inside
# defined names (e.g. the hidden _FilterDatabase name) to use the new
# sheet name when the active sheet is renamed.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Name = "البيانات"

# The "الصنف" defined name pointed at a broken reference on the old sheet
# ('مثال على الاخراج'!#REF!). Re-point it explicitly at the renamed sheet
# so it reads "البيانات!#REF!" after the rename, matching how Excel keeps
# the sheet qualifier on broken references belonging to the renamed sheet.
foreach ($n in $wb.Names) {
    if ($n.Name -eq "الصنف") {
        $n.RefersTo = "=البيانات!#REF!"
    }
}

# Move/save the active selection to A12 on the sheet.
$ws.Range("A12").Select()
